# daily auto push: 2026-02-21 22:39 UTC
# Insert a new data point for 2026/02/22 (日) at 05:00, right after the
# existing 2026/02/22 (日) 02:00 row (row 831). This pushes every row
# below it (old 831..872) down by one (new 832..873), which matches the
# dimension growing from A1:D872 to A1:D873.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 831 (old row 831 and everything
# below shifts down to 832, 833, ...).
$ws.Rows(831).Insert()

# Column A holds a date-looking string ("2026/02/22") that must stay
# literal text (matching every other cell in the column) instead of
# being auto-parsed into a date serial number by the COM layer. Forcing
# the cell to Text format before the assignment, then clearing the
# format afterwards, keeps the stored value as plain text while leaving
# the cell's style back at the sheet default (same as its neighbours).
$ws.Range("A831").NumberFormat = "@"
$ws.Range("A831").Value = "2026/02/22"
$ws.Range("A831").ClearFormats()

$ws.Range("B831").Value = "日"
$ws.Range("C831").Value = 5
$ws.Range("D831").Value = 42
